$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "96.467.81"
Set-TextValue "E2" "  -0.32%  "
Set-TextValue "D3" "3.668.56"
Set-TextValue "E3" "  +2.53%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "241.42"
Set-TextValue "E5" "  -0.04%  "
Set-TextValue "E6" "  +11.97%  "
Set-TextValue "D7" "666.77"
Set-TextValue "E7" "  +2.02%  "
Set-TextValue "D8" "0.425"
Set-TextValue "E8" "  +3.25%  "
Set-TextValue "E9" "  +0.07%  "
Set-TextValue "E10" "  +0.06%  "
Set-TextValue "D11" "3.667.20"
Set-TextValue "E11" "  +2.49%  "
Set-TextValue "D12" "44.85"
Set-TextValue "E12" "  +2.87%  "
Set-TextValue "E13" "  +0.51%  "
Set-TextValue "D14" "6.64"
Set-TextValue "E14" "  +3.71%  "
Set-TextValue "D15" "4.350.05"
Set-TextValue "E15" "  +2.55%  "
Set-TextValue "D17" "96.371.78"
Set-TextValue "E17" "  -0.24%  "
Set-TextValue "D18" "8.87"
Set-TextValue "E18" "  +14.08%  "
Set-TextValue "D19" "3.668.69"
Set-TextValue "E19" "  +2.91%  "
Set-TextValue "D20" "12.72"
Set-TextValue "E20" "  +0.63%  "
Set-TextValue "D21" "18.32"
Set-TextValue "E21" "  +1.73%  "
Set-TextValue "D22" "0.533"
Set-TextValue "E22" "  -2.08%  "
Set-TextValue "D23" "524.14"
Set-TextValue "E23" "  +3.02%  "
Set-TextValue "D24" "3.44"
Set-TextValue "E24" "  +1.10%  "
Set-TextValue "D25" "0.0000204"
Set-TextValue "E25" "  +1.54%  "
Set-TextValue "D26" "6.91"
Set-TextValue "E26" "  -0.81%  "
Set-TextValue "D27" "102.54"
Set-TextValue "E27" "  +5.90%  "
Set-TextValue "E28" "  -1.14%  "
Set-TextValue "B29" "Hedera"
Set-TextValue "C29" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D29" "0.166"
Set-TextValue "E29" "  +7.27%  "
Set-TextValue "B30" "PancakeSwap"
Set-TextValue "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "3.03"
Set-TextValue "E30" "  -0.04%  "
Set-TextValue "B31" "InternetComputer(DFINITY)"
Set-TextValue "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "12.11"
Set-TextValue "E31" "  +5.13%  "
Set-TextValue "B32" "Dai"
Set-TextValue "C32" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D32" "0.996"
Set-TextValue "E32" "  -0.37%  "
Set-TextValue "B33" "Cronos"
Set-TextValue "C33" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D33" "0.185"
Set-TextValue "E33" "  +0.52%  "
Set-TextValue "B34" "Fetch.AI"
Set-TextValue "C34" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D34" "1.81"
Set-TextValue "E34" "  +10.61%  "
Set-TextValue "B35" "EthereumClassic"
Set-TextValue "C35" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "32.67"
Set-TextValue "E35" "  +3.97%  "
Set-TextValue "B36" "Binance-PegBSC-USD"
Set-TextValue "C36" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  +0.19%  "
Set-TextValue "B37" "PolygonEcosystemToken"
Set-TextValue "C37" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D37" "0.585"
Set-TextValue "E37" "  +2.57%  "
Set-TextValue "B38" "Bittensor"
Set-TextValue "C38" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D38" "617.38"
Set-TextValue "E38" "  -1.36%  "
Set-TextValue "B39" "RenderToken"
Set-TextValue "C39" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D39" "8.73"
Set-TextValue "E39" "  -1.55%  "
Set-TextValue "B40" "EnergySwap"
Set-TextValue "C40" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D40" "42.81"
Set-TextValue "E40" "  +30.28%  "
Set-TextValue "B41" "Kaspa"
Set-TextValue "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.159"
Set-TextValue "E41" "  +4.92%  "
Set-TextValue "B42" "ARBITRUM"
Set-TextValue "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.950"
Set-TextValue "E42" "  +4.69%  "
Set-TextValue "B43" "ImmutableX"
Set-TextValue "C43" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D43" "1.94"
Set-TextValue "E43" "  +6.17%  "
Set-TextValue "B44" "USDe"
Set-TextValue "C44" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  +0.00%  "
Set-TextValue "B45" "Filecoin"
Set-TextValue "C45" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "6.16"
Set-TextValue "E45" "  +6.33%  "
Set-TextValue "B46" "VeChain"
Set-TextValue "C46" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0454"
Set-TextValue "E46" "  +5.42%  "
Set-TextValue "B47" "Algorand"
Set-TextValue "C47" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D47" "0.426"
Set-TextValue "E47" "  +19.10%  "
Set-TextValue "B48" "Stacks"
Set-TextValue "C48" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D48" "2.29"
Set-TextValue "E48" "  -1.03%  "
Set-TextValue "B49" "WhiteBITCoin"
Set-TextValue "C49" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D49" "23.60"
Set-TextValue "E49" "  +0.23%  "
Set-TextValue "B50" "Cosmos"
Set-TextValue "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "8.49"
Set-TextValue "E50" "  +1.81%  "
Set-TextValue "B51" "MantraDAO"
Set-TextValue "C51" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D51" "3.55"
Set-TextValue "E51" "  +0.86%  "
